$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6258
$ws.Range("C21").Value = 990
$ws.Range("D21").Value = 5671749
$ws.Range("E21").Value = 906.3197507190796
$ws.Range("F21").Value = 8.626974483596594
$ws.Range("G21").Value = 4.430379746835444
$ws.Range("H21").Value = 29.46325609933138
